$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" column (G) for the
# 1d4f27e1-7e39-43db-9fd3-86b508e050fa.md row is reported twice
# (rows 2 and 3 share the same shared-string value).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 20:16:39"
$wsOverview.Range("G3").Value = "2016-08-24 20:16:39"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E): "ht" -> "mt"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-08-24 20:16:34"
$wsZhCn.Range("H3").Value = "2016-08-24 20:16:34"
# Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-08-24 20:16:52"
$wsZhCn.Range("K3").Value = "2016-08-24 20:16:52"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority column (E): "ht" -> "mt" (shares the same shared-string value
# as the zh-cn sheet's E2/E3 cells)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
# Correspond Handoff Datetime column (H) - shares the same original value
# as the Overview sheet's G2/G3 cells
$wsDeDe.Range("H2").Value = "2016-08-24 20:16:39"
$wsDeDe.Range("H3").Value = "2016-08-24 20:16:39"
# Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-08-24 20:16:59"
$wsDeDe.Range("K3").Value = "2016-08-24 20:16:59"
